# Updated cryptos list - apply price and volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format first so numeric-looking strings
# (e.g. "246.22") are stored as text, matching the original inlineStr cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "25.980.67"
$ws.Range("E2").Value = "  +0.60%  "

$ws.Range("D3").Value = "1.740.12"
$ws.Range("E3").Value = "  +0.38%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").Value = "246.22"
$ws.Range("E5").Value = "  +4.03%  "

$ws.Range("D7").Value = "0.5021"
$ws.Range("E7").Value = "  -1.75%  "

$ws.Range("D8").Value = "0.2736"
$ws.Range("E8").Value = "  +1.01%  "

$ws.Range("D9").Value = "0.06180"
$ws.Range("E9").Value = "  +1.62%  "

$ws.Range("D10").Value = "1.746.12"
$ws.Range("E10").Value = "  +0.80%  "

$ws.Range("D11").Value = "0.07251"
$ws.Range("E11").Value = "  +1.49%  "

$ws.Range("D12").Value = "0.6529"
$ws.Range("E12").Value = "  +3.45%  "

$ws.Range("D13").Value = "15.08"
$ws.Range("E13").Value = "  +1.30%  "

$ws.Range("D14").Value = "4.678"
$ws.Range("E14").Value = "  +2.21%  "

$ws.Range("D15").Value = "77.56"
$ws.Range("E15").Value = "  +0.95%  "

$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  -0.07%  "

$ws.Range("E17").Value = "  -0.10%  "

$ws.Range("D18").Value = "26.012.32"
$ws.Range("E18").Value = "  +0.67%  "

$ws.Range("D19").Value = "11.92"
$ws.Range("E19").Value = "  +2.19%  "

$ws.Range("D20").Value = "0.000006853"
$ws.Range("E20").Value = "  +2.68%  "

$ws.Range("D21").Value = "1.968.84"
$ws.Range("E21").Value = "  +0.46%  "

$ws.Range("D22").Value = "4.496"
$ws.Range("E22").Value = "  +5.51%  "

$ws.Range("D23").Value = "8.697"
$ws.Range("E23").Value = "  +1.17%  "

$ws.Range("D24").Value = "5.395"
$ws.Range("E24").Value = "  +3.65%  "

$ws.Range("D25").Value = "135.51"
$ws.Range("E25").Value = "  -2.82%  "

$ws.Range("D26").Value = "1.513"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("D27").Value = "15.27"
$ws.Range("E27").Value = "  +0.95%  "

$ws.Range("D28").Value = "1.782"
$ws.Range("E28").Value = "  +1.70%  "

$ws.Range("D29").Value = "105.45"
$ws.Range("E29").Value = "  +0.39%  "

$ws.Range("D30").Value = "3.950"
$ws.Range("E30").Value = "  +1.88%  "

$ws.Range("D31").Value = "0.08149"
$ws.Range("E31").Value = "  -2.12%  "

$ws.Range("D32").Value = "3.673"
$ws.Range("E32").Value = "  +2.95%  "

$ws.Range("D33").Value = "0.04698"
$ws.Range("E33").Value = "  +3.36%  "

$ws.Range("E34").Value = "  +1.44%  "

$ws.Range("D35").Value = "0.9933"
$ws.Range("E35").Value = "  +1.55%  "

$ws.Range("D36").Value = "0.6095"
$ws.Range("E36").Value = "  -1.34%  "

$ws.Range("D37").Value = "2.764"
$ws.Range("E37").Value = "  +2.67%  "

$ws.Range("D38").Value = "0.01618"
$ws.Range("E38").Value = "  +1.98%  "

$ws.Range("D39").Value = "1.923"
$ws.Range("E39").Value = "  +0.71%  "

$ws.Range("E40").Value = "  +0.01%  "

$ws.Range("D41").Value = "100.99"
$ws.Range("E41").Value = "  +4.03%  "

$ws.Range("D42").Value = "0.7945"
$ws.Range("E42").Value = "  +8.43%  "

$ws.Range("D43").Value = "0.3897"
$ws.Range("E43").Value = "  +1.92%  "

$ws.Range("D44").Value = "5.012"
$ws.Range("E44").Value = "  +1.85%  "

$ws.Range("D45").Value = "0.1164"
$ws.Range("E45").Value = "  +3.29%  "

$ws.Range("D46").Value = "6.308"
$ws.Range("E46").Value = "  +2.39%  "

$ws.Range("D47").Value = "55.61"
$ws.Range("E47").Value = "  +2.23%  "

$ws.Range("D48").Value = "0.05292"
$ws.Range("E48").Value = "  +0.34%  "

$ws.Range("D49").Value = "30.80"
$ws.Range("E49").Value = "  +1.40%  "

$ws.Range("D50").Value = "7.631"
$ws.Range("E50").Value = "  +0.87%  "

$ws.Range("D51").Value = "0.3465"
$ws.Range("E51").Value = "  +1.91%  "

# Restore the default "Normal" style on column D so no stray number-format
# style is left applied to the cells (keeps formatting identical to original).
$ws.Range("D2:D51").Style = "Normal"
